$wb = $excel.ActiveWorkbook

# 1) Status text: "Ready for handoff" -> "In Translation"
#    (appears on Overview!E2, Overview!F2, zh-cn!C2, de-de!C2 - all share the
#    same underlying string, so a workbook-wide replace keeps them in sync)
foreach ($ws in $wb.Worksheets) {
    $ws.Cells.Replace("Ready for handoff", "In Translation")
}

# 2) Narrow the per-language status columns
#    Overview: columns E (zh-cn) and F (de-de)
#    zh-cn / de-de detail sheets: column C (Status)
$newWidth = 12.5

$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Columns.Item(5).ColumnWidth = $newWidth
$wsOverview.Columns.Item(6).ColumnWidth = $newWidth

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Columns.Item(3).ColumnWidth = $newWidth

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Columns.Item(3).ColumnWidth = $newWidth
